# Improve recording UI — add a new "Generate alpha (Blend+Add)" localization
# entry right after the existing "GenerateAlpha" row, shifting every row
# below it down by one (Excel-style row insert, which preserves styles,
# row heights and all the other unrelated rows/cells untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 34 (pushes old row 34.."Height".. down to 35, etc.)
$ws.Rows("34:34").Insert()

# Fill in the new localization row with the new key/en/ja strings.
# (Rows("34:34").Insert() already copied the surrounding row's formatting
# onto the new row, so no separate style assignment is needed here.)
$ws.Range("A34").Value = "GenerateAlpha2"
$ws.Range("B34").Value = "Generate alpha(Blend+Add)"
$ws.Range("C34").Value = "生成(ブレンド+加算)"

# Reflect the author's final navigation state in the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("C35").Select()
